$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp in header cell A1
$ws.Range("A1").Value = "Datos actualizados a 7 de Abril de 2020 a las 15:52"

# Germany (Alemania) - row 7
$ws.Range("B7").Value = 104199
$ws.Range("C7").Value = 824
$ws.Range("E7").Value = 66276
$ws.Range("G7").Value = 32
$ws.Range("H7").Value = 1842

# Rumania - row 32
$ws.Range("E32").Value = 3760
$ws.Range("G32").Value = 21
$ws.Range("H32").Value = 197

# Pakistan - row 33
$ws.Range("B33").Value = 4005
$ws.Range("C33").Value = 239
$ws.Range("E33").Value = 3521

# Islandia moves above Colombia (row 53 becomes Islandia, row 54 becomes Colombia)
$ws.Range("A53").Value = "Islandia"
$ws.Range("B53").Value = 1586
$ws.Range("C53").Value = 24
$ws.Range("D53").Value = 559
$ws.Range("E53").Value = 1021
$ws.Range("F53").Value = 11
$ws.Range("H53").Value = 6

$ws.Range("A54").Value = "Colombia"
$ws.Range("B54").Value = 1579
$ws.Range("C54").Value = 0
$ws.Range("D54").Value = 88
$ws.Range("E54").Value = 1445
$ws.Range("F54").Value = 50
$ws.Range("H54").Value = 46

# Eslovenia - row 63
$ws.Range("F63").Value = 30

# Eslovaquia - row 80
$ws.Range("D80").Value = 13
$ws.Range("E80").Value = 566

# Islas Feroe - row 113
$ws.Range("D113").Value = 120
$ws.Range("E113").Value = 64

# Zambia - row 148
$ws.Range("D148").Value = 7
$ws.Range("E148").Value = 31
